# Insert a new data row at row 486 (shifts existing rows 486:549 down to 487:550)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(486).Insert()

# Populate the newly inserted row 486 with the new weekly price entry.
$ws.Cells.Item(486, 1).Value = 11
$ws.Cells.Item(486, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(486, 3).Value = "Bíobío"
$ws.Cells.Item(486, 4).Value = 45154
$ws.Cells.Item(486, 5).Value = 8
$ws.Cells.Item(486, 6).Value = 100114014
$ws.Cells.Item(486, 7).Value = "Betarraga"
$ws.Cells.Item(486, 8).Value = "Sin especificar"
$ws.Cells.Item(486, 9).Value = "Primera"
$ws.Cells.Item(486, 10).Value = 600
$ws.Cells.Item(486, 11).Value = 600
$ws.Cells.Item(486, 12).Value = 650
$ws.Cells.Item(486, 13).Value = 625
$ws.Cells.Item(486, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(486, 15).Value = "Región Metropolitana"
$ws.Cells.Item(486, 16).Value = 125
$ws.Cells.Item(486, 17).Value = 5
$ws.Cells.Item(486, 18).Value = "Hortaliza"
